# Daily update at 8 AM UTC
# Appends a new row of data (row 83) to the bottom of the tracking sheet,
# and fixes up the "last row" date style so only the new last row keeps it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last data row before this update is 82; new row goes to 83.
$lastRow = 82
$newRow = $lastRow + 1

# Remember the special "final row" date format currently applied to A82.
$finalRowFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

# Reset the previous last row's date cell (A82) back to the regular format,
# matching the format used by every other non-final row (e.g. A81).
$ws.Cells.Item($lastRow, 1).NumberFormat = $ws.Cells.Item($lastRow - 1, 1).NumberFormat

# Add the new day's data.
$ws.Cells.Item($newRow, 1).Value = 45823
$ws.Cells.Item($newRow, 2).Value = 355
$ws.Cells.Item($newRow, 3).Value = 355
$ws.Cells.Item($newRow, 4).Value = 360

# The new last row (A83) takes on the "final row" date format that A82 used to have.
$ws.Cells.Item($newRow, 1).NumberFormat = $finalRowFormat
